# Expanded to allow for off the money pricing example.
#
# Adds five new pricing-example worksheets ("0.5x10.15" .. "0.5x10.19") after
# the existing "0.5x10.14" sheet, each a duplicate of that sheet's layout
# (dates/nominal/fwd/strike/expiry/rf/vol/tenor/type/price/interpolated
# key-value table). The fourth of the new sheets ("0.5x10.18") is an
# off-the-money example and gets a different strike, vol and resulting
# price; the others keep the original at-the-money values.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("0.5x10.14")

# Duplicate the template sheet (this also copies its formatting/styles and
# page setup exactly), placing each copy immediately after its source so the
# new tabs land at the end of the workbook in order.
$src.Copy($null, $src)
$sheet15 = $wb.Worksheets.Item($src.Index + 1)
$sheet15.Name = "0.5x10.15"

$sheet15.Copy($null, $sheet15)
$sheet16 = $wb.Worksheets.Item($sheet15.Index + 1)
$sheet16.Name = "0.5x10.16"

$sheet16.Copy($null, $sheet16)
$sheet17 = $wb.Worksheets.Item($sheet16.Index + 1)
$sheet17.Name = "0.5x10.17"

$sheet17.Copy($null, $sheet17)
$sheet18 = $wb.Worksheets.Item($sheet17.Index + 1)
$sheet18.Name = "0.5x10.18"

$sheet18.Copy($null, $sheet18)
$sheet19 = $wb.Worksheets.Item($sheet18.Index + 1)
$sheet19.Name = "0.5x10.19"

# "0.5x10.18" is the off-the-money pricing example: different strike/vol,
# and the resulting recalculated price.
$sheet18.Range("B5").Value = 2.2
$sheet18.Range("B8").Value = 50.63
$sheet18.Range("B11").Value = 0.8021397878237978
